$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.788.42"
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").Value = "3.422.29"
$ws.Range("E3").Value = "  -2.62%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.79"
$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.58"
$ws.Range("E6").Value = "  -4.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.422.64"
$ws.Range("E8").Value = "  -2.58%  "

$ws.Range("E9").Value = "  -5.91%  "

$ws.Range("E10").Value = "  -9.79%  "

$ws.Range("E11").Value = "  -9.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("E12").Value = "  -6.98%  "

$ws.Range("D13").Value = "3.997.71"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  -9.13%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.41"
$ws.Range("E15").Value = "  -6.76%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.416.08"
$ws.Range("E16").Value = "  -2.56%  "

$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Value = "64.620.86"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.56"
$ws.Range("E19").Value = "  -11.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.80"
$ws.Range("E20").Value = "  -5.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  -5.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.83"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  -7.79%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.09"
$ws.Range("E25").Value = "  -6.35%  "

$ws.Range("D26").Value = "3.552.64"
$ws.Range("E26").Value = "  -2.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  -7.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  -6.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("E30").Value = "  -9.26%  "

$ws.Range("E31").Value = "  -10.79%  "

$ws.Range("D32").Value = "3.427.69"
$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.93"
$ws.Range("E34").Value = "  -5.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.141"
$ws.Range("E35").Value = "  -8.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.04"
$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  -11.30%  "

$ws.Range("E38").Value = "  -10.69%  "

$ws.Range("E39").Value = "  -7.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.63"
$ws.Range("E40").Value = "  -11.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0757"
$ws.Range("E41").Value = "  -6.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.809"
$ws.Range("E42").Value = "  -5.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.80"
$ws.Range("E44").Value = "  -7.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.30"
$ws.Range("E45").Value = "  -13.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("E46").Value = "  -8.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +4.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.61"
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.46"
$ws.Range("E49").Value = "  -8.17%  "

$ws.Range("E50").Value = "  -13.23%  "

$ws.Range("D51").Value = "2.165.57"
$ws.Range("E51").Value = "  -6.86%  "
